# 편의시설 영업시간.xlsx - fix two typo'd operating-hour strings and tidy up
# the leftover "apply alignment" formatting on the location column (A2:A15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "11:00~!7:30" (stray "!") -> "11:00~7:30" for 츄밥
$ws.Range("C12").Value = "11:00~7:30"

# Fix "10:00~17:0" (missing trailing 0) -> "10:00~17:00" for 학생 누리관 / 에땅 샌드위치
$ws.Range("C8").Value = "10:00~17:00"

# The A column had an extra/duplicate "apply alignment" style left over on
# rows 2-15; clear it back to the sheet's default style.
$ws.Range("A2:A15").ClearFormats()

# Leave the selection where the editor ended up after making the last fix.
$ws.Range("C8").Select() | Out-Null
